$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "VOLTAREN SR 100MG 20 F.C.TAB." (item #16, row 22)
# This shifts rows 23-27 up by one.
$ws.Rows("22:22").Delete()

# The item-number column (A) keeps a simple sequential numbering for the
# remaining items, so renumber the rows that shifted up.
$ws.Range("A22").Value = 16
$ws.Range("A23").Value = 17
$ws.Range("A24").Value = 18

# Update the timestamp string (now located in A26 after the shift)
$ws.Range("A26").Value = "Tuesday, 10 June, 2025 10:51 AM"

# Update the grand total (sum of "sell price" column) to reflect the removed item
$ws.Range("P25").Value = 1008.58
